$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $found = $r.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $found) {
        throw "Find failed for: $find"
    }
}

# 1. Date update
Replace-Text "2024-10-31" "2024-11-01"

# 2. "higher level of the trait" -> "higher trait level"
Replace-Text "indicating which stimulus is perceived to possess a higher level of the trait." "indicating which stimulus is perceived to possess a higher trait level."

# 3. "CJ has been successfully employed" -> "This method has been successfully employed"
Replace-Text "underlying trait of interest. CJ has been successfully employed in assessing the quality of written texts" "underlying trait of interest. This method has been successfully employed in assessing the quality of written texts"

# 4. "On the other hand" -> "Additionally" and "traits being measured" -> "traits under measurement"
Replace-Text ". On the other hand, research on validity suggests that scores generated by CJ can accurately represent the traits being measured" ". Additionally, research on validity suggests that scores generated by CJ can accurately represent the traits under measurement"

# 5. "Finally, research" -> ", while research"
Replace-Text ". Finally, research on practical applicability highlights the method" ", while research on practical applicability highlights the method"

# 6. "growing number of studies on CJ" -> "growing number of CJ studies"
Replace-Text "Nevertheless, despite the growing number of studies on CJ, unsystematic and fragmented" "Nevertheless, despite the growing number of CJ studies, unsystematic and fragmented"

# 7. "these scores have been used" -> "the scores have been used"
Replace-Text "The literature indicates that these scores have been used to identify" "The literature indicates that the scores have been used to identify"

# 8. "other scoring methods" -> "other assessment methods"
Replace-Text ", calculate correlations with other scoring methods" ", calculate correlations with other assessment methods"

# 9. "latent trait of interest" -> "underlying trait of interest"
Replace-Text ", or test hypotheses related to the latent trait of interest" ", or test hypotheses related to the underlying trait of interest"

# 10. "conducting additional analyses and tests can inflate" -> "conducting separate analyses and tests can inflate"
Replace-Text "Ignoring this uncertainty when conducting additional analyses and tests can inflate" "Ignoring this uncertainty when conducting separate analyses and tests can inflate"

# 11. Final paragraph rewrite - split run into multiple runs with a new citation
Replace-Text "recommend conducting these analyses and tests within a structural model that accounts for both the scores and their uncertainties, rather than treating them separately. Thus, an integrated approach combining CJ" "recommend conducting these analyses and tests within a structural model. A structural model specify how different manifest or latent variables influence the latent variable of interest (Everitt and Skrondal 2010). This approach allows analyses that can account for both the scores and their uncertainties simultaneously, rather than treating them as separate elements. Therefore, an integrated approach that combines CJ"
